$d = $word.ActiveDocument

$replacements = @(
    @{Old = "Un monde de saveurs dans une tasse"; New = "Thé Chai : Un monde de saveur dans une tasse"},
    @{Old = "L’alliance parfaite de la santé et du plaisir"; New = "Thé Chai : Le mélange parfait de santé et de plaisir"},
    @{Old = "Plus qu’un simple thé, un mode de vie"; New = "Thé Chai : Plus que le thé, un mode de vie"},
    @{Old = "Une boisson pour toutes les raisons, en toute saison"; New = "Thé Chai : Une boisson pour toutes les saisons et les raisons"},
    @{Old = "L’ultime plaisir des sens"; New = "Thé Chai : L’indulgence ultime pour vos sens"},
    @{Old = "Une douce évasion du quotidien"; New = "Thé Chai : Une évasion douce du quotidien"},
    @{Old = "Partager la chaleur, partager l’amour"; New = "Thé Chai : Partager la chaleur, partager l’amour"}
)

foreach ($r in $replacements) {
    $rng = $d.Content
    $rng.Find.Execute($r.Old, $true, $false, $false, $false, $false, $true, 1, $false, $r.New, 2)
}
